# edit.ps1 - applies the Screw BOSC 2017 manuscript edit:
#   1) Splits long single-run sentences into one run per sentence
#      (with a separate single-space run between sentences), in the
#      first three body paragraphs.
#   2) Fixes wording: "especially in contexts where" ->
#      "especially in rare cell types and contexts where"
#   3) Fixes wording: "enable complete reproducible" ->
#      "enable completely reproducible"
#   4) Changes the <w:nsid> of the abstractNum with abstractNumId=990
#      from aad61e29 to cb1bb5ec.

$d = $word.ActiveDocument

function Replace-ParagraphRuns($Paragraph, $Texts) {
    $pStart = $Paragraph.Range.Start
    $pEnd = $Paragraph.Range.End

    # Paragraph.Range.End includes the trailing paragraph mark; stop one
    # character short of it so only the paragraph's content is replaced.
    $target = $d.Range($pStart, $pEnd - 1)

    $runsXml = ""
    foreach ($t in $Texts) {
        $escaped = $t.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        $runsXml += '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
    }

    $pPrXml = ""
    $style = $null
    try { $style = $Paragraph.Style.NameLocal } catch { $style = $null }
    if ($style) {
        $pPrXml = '<w:pPr><w:pStyle w:val="' + $style + '"/></w:pPr>'
    }

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $pPrXml + $runsXml + '</w:p>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)
}

# --- Paragraph 2 (FirstParagraph style): intro paragraph ---------------
$p2 = $d.Paragraphs.Item(2)
$texts2 = @(
    "DNA methylation is a heritable epigenetic mark that shows a strong correlation with transcriptional activity.",
    " ",
    "The gold standard for detecting DNA methylation is whole genome bisulfite sequencing (WGBS).",
    " ",
    "Recently, WGBS has been performed successfully on single cells (SC-WGBS)",
    " ",
    "[1]",
    ".",
    " ",
    "The resulting data represents a fundamental shift in the capacity to measure and interpret DNA methylation, especially in rare cell types and contexts where subtle cell-to-cell heterogeneity is crucial, such as in stem cells or cancer.",
    " ",
    "However, SC-WGBS comes with unique technical challenges which require new analysis techniques to address.",
    " ",
    "Furthermore, although some tools have been published, and several existing studies have tended to use similar methods, no standardized pipeline for the analysis of SC-WGBS yet exists."
)
Replace-ParagraphRuns $p2 $texts2

# --- Paragraph 3 (BodyText style): reproducibility paragraph -----------
$p3 = $d.Paragraphs.Item(3)
$texts3 = @(
    "Simultaneously, there has been a drive within bioinformatics towards improved reproducibility.",
    " ",
    "Textual descriptions of bioinformatic analyses are deeply inadequate, and often require `"forensic bioinformatics`" to reproduce",
    " ",
    "[2]",
    ".",
    " ",
    "Exact code, accompanied by exact software versions used, is needed to recreate the exact results of a study.",
    " ",
    "Common Workflow Language (CWL) provides a framework for specifying complete workflows, while Docker allows for bundling of the exact software used in an analysis within a container that can be executed anywhere.",
    " ",
    "Together, these have the potential, via repositories such as Dockstore",
    " ",
    "[3]",
    ", to enable completely reproducible bioinformatics research."
)
Replace-ParagraphRuns $p3 $texts3

# --- Paragraph 4 (BodyText style): Screw presentation paragraph --------
$p4 = $d.Paragraphs.Item(4)
$texts4 = @(
    "Here we present Screw (Single Cell Reproducible Epigenomics Workfow).",
    " ",
    "Screw is a collection of standard tools and workflows for analysing SC-WGBS data, implemented in CWL, and with an accompanying Docker image.",
    " ",
    "Screw is intended to provide the parts to build fully-reproducible SC-WGBS analyses.",
    " ",
    "Tools provided include quality control visualization, clustering and visualisation of cells by pairwise dissimilarity measures, construction of recapitulated-bulk methylomes from single cells of the same lineage, generation of bigWig methylation tracks for downstream visualization, and wrappers around published tools such as DeepCpG",
    " ",
    "[4]",
    " ",
    "and LOLA",
    " ",
    "[5]",
    ".",
    " ",
    "Screw has the added benefit that CWL's compatibility with interactive GUI-based workflow tools such as Galaxy can lower the barriers to use for less-technical wet lab biologist users."
)
Replace-ParagraphRuns $p4 $texts4

# --- numbering.xml: change nsid for abstractNumId=990 -------------------
$found = $d.Content.Find.Execute("aad61e29", $true, $false, $false, $false, $false, $true, 1, $false, "cb1bb5ec", 2)
Write-Host "nsid replace found:" $found
